$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Abstract" column (F), shifting the
# old F ("Abstract" text) -> G and old G (empty, s=2) -> H.
$ws.Columns.Item(6).Insert()

# The Insert() shifts the trailing "whole sheet" <col> range (…:16384) by one
# to …:16385; drop the now-phantom last column to restore the 1..16384 range.
$ws.Columns.Item(16384).Delete()

# Fill the new column F with a simple 0-based row index (rows 2..51 -> 0..49).
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 6).Value = $r - 2
}

# Narrow the new column F and widen the (shifted) column G slightly.
$ws.Columns.Item(6).ColumnWidth = 3.33
$ws.Columns.Item(7).ColumnWidth = 10.8

# Move the active selection to B1.
$ws.Range("B1").Select()
